$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2005.1765
$ws.Range("I125").Value = 866
$ws.Range("J125").Value = 2157.0667
$ws.Range("K125").Value = 7794
$ws.Range("L125").Value = 19413.6003
$ws.Range("M125").Value = -5334
$ws.Range("N125").Value = -24333.6003
$ws.Range("H138").Value = 1269.5416
$ws.Range("I138").Value = 1015.9778
$ws.Range("J138").Value = 5073
$ws.Range("K138").Value = 3047.9334
$ws.Range("L138").Value = 15219
$ws.Range("M138").Value = 2092.0666
$ws.Range("N138").Value = -25499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 2002065.2
$ws.Range("I36").Value = 2002065.2
$ws.Range("K36").Value = 2002065.2
$ws.Range("M36").Value = -2001719.2
$ws.Range("H61").Value = 1519.88
$ws.Range("I61").Value = 1361.6774
$ws.Range("J61").Value = 1778
$ws.Range("K61").Value = 1361.6774
$ws.Range("L61").Value = 1778
$ws.Range("M61").Value = -1149.6774
$ws.Range("N61").Value = -2202
$ws.Range("H88").Value = 2275.8
$ws.Range("I88").Value = 1299.6666
$ws.Range("J88").Value = 3740
$ws.Range("K88").Value = 1299.6666
$ws.Range("L88").Value = 3740
$ws.Range("M88").Value = -893.6666
$ws.Range("N88").Value = -4552
$ws.Range("H91").Value = 2275.8
$ws.Range("I91").Value = 1299.6666
$ws.Range("J91").Value = 3740
$ws.Range("K91").Value = 1299.6666
$ws.Range("L91").Value = 3740
$ws.Range("M91").Value = 104.3334
$ws.Range("N91").Value = -6548
$ws.Range("H136").Value = 1519.88
$ws.Range("I136").Value = 1361.6774
$ws.Range("J136").Value = 1778
$ws.Range("K136").Value = 4085.0322
$ws.Range("L136").Value = 5334
$ws.Range("M136").Value = -1535.0322
$ws.Range("N136").Value = -10434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 13000
$ws.Range("I33").Value = 13000
$ws.Range("K33").Value = 13000
$ws.Range("M33").Value = -12664
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41372
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126864
$ws.Range("H86").Value = 1693.125
$ws.Range("I86").Value = 1488.75
$ws.Range("J86").Value = 1897.5
$ws.Range("K86").Value = 1488.75
$ws.Range("L86").Value = 1897.5
$ws.Range("M86").Value = -365.75
$ws.Range("N86").Value = -4143.5
$ws.Range("H89").Value = 1693.125
$ws.Range("I89").Value = 1488.75
$ws.Range("J89").Value = 1897.5
$ws.Range("K89").Value = 7443.75
$ws.Range("L89").Value = 9487.5
$ws.Range("M89").Value = -1827.75
$ws.Range("N89").Value = -20719.5
$ws.Range("H134").Value = 6154.385
$ws.Range("I134").Value = 3434.9167
$ws.Range("J134").Value = 7363.037
$ws.Range("K134").Value = 10304.7501
$ws.Range("L134").Value = 22089.111
$ws.Range("M134").Value = -7769.750100000001
$ws.Range("N134").Value = -27159.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H58").Value = 3075.988
$ws.Range("I58").Value = 1822.2458
$ws.Range("J58").Value = 6401.1304
$ws.Range("K58").Value = 1822.2458
$ws.Range("L58").Value = 6401.1304
$ws.Range("M58").Value = -1619.2458
$ws.Range("N58").Value = -6807.1304
$ws.Range("H99").Value = 6000
$ws.Range("I99").Value = 6000
$ws.Range("K99").Value = 6000
$ws.Range("M99").Value = -4502
$ws.Range("H110").Value = 39333.332
$ws.Range("I110").Value = 40000
$ws.Range("J110").Value = 39000
$ws.Range("K110").Value = 40000
$ws.Range("L110").Value = 39000
$ws.Range("M110").Value = -35910
$ws.Range("N110").Value = -47180
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 6000
$ws.Range("K126").Value = 18000
$ws.Range("M126").Value = -15530
$ws.Range("H136").Value = 3075.988
$ws.Range("I136").Value = 1822.2458
$ws.Range("J136").Value = 6401.1304
$ws.Range("K136").Value = 5466.7374
$ws.Range("L136").Value = 19203.3912
$ws.Range("M136").Value = -2916.7374
$ws.Range("N136").Value = -24303.3912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1996
$ws.Range("J41").Value = 2021.3684
$ws.Range("L41").Value = 6064.1052
$ws.Range("N41").Value = -6740.1052
$ws.Range("H137").Value = 3503.239
$ws.Range("I137").Value = 2278.625
$ws.Range("J137").Value = 3761.0527
$ws.Range("K137").Value = 6835.875
$ws.Range("L137").Value = 11283.1581
$ws.Range("M137").Value = -1735.875
$ws.Range("N137").Value = -21483.1581
$ws.Range("H138").Value = 2281.3928
$ws.Range("I138").Value = 1554.1666
$ws.Range("K138").Value = 4662.4998
$ws.Range("M138").Value = 477.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5937.7334
$ws.Range("I70").Value = 5263.533
$ws.Range("J70").Value = 6611.933
$ws.Range("K70").Value = 5263.533
$ws.Range("L70").Value = 6611.933
$ws.Range("M70").Value = -4993.533
$ws.Range("N70").Value = -7151.933
$ws.Range("H73").Value = 5937.7334
$ws.Range("I73").Value = 5263.533
$ws.Range("J73").Value = 6611.933
$ws.Range("K73").Value = 5263.533
$ws.Range("L73").Value = 6611.933
$ws.Range("M73").Value = -4327.533
$ws.Range("N73").Value = -8483.933000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1991.0769
$ws.Range("I7").Value = 1907
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1907
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1795
$ws.Range("N7").Value = -3224
$ws.Range("H126").Value = 1991.0769
$ws.Range("I126").Value = 1907
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5721
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3251
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 37041730
$ws.Range("I132").Value = 66672604
$ws.Range("K132").Value = 200017812
$ws.Range("M132").Value = -200015282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 16625
$ws.Range("J94").Value = 16625
$ws.Range("L94").Value = 16625
$ws.Range("N94").Value = -18427
$ws.Range("H126").Value = 2193.2778
$ws.Range("I126").Value = 3054.9092
$ws.Range("J126").Value = 839.2857
$ws.Range("K126").Value = 9164.7276
$ws.Range("L126").Value = 2517.8571
$ws.Range("M126").Value = -6694.7276
$ws.Range("N126").Value = -7457.8571
